# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price (D) and Volume(1h) (E) text cells are updated per-row; values that
# look numeric are apostrophe-prefixed so Excel keeps them as text (preserving
# trailing zeros / exact digit strings instead of coercing to Number).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.458.95"
$ws.Range("E2").Value = "  -3.80%  "
$ws.Range("D3").Value = "2.511.95"
$ws.Range("E3").Value = "  -4.53%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'576.41"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").Value = "'166.02"
$ws.Range("E6").Value = "  -4.45%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.515"
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("D9").Value = "2.510.05"
$ws.Range("E9").Value = "  -4.55%  "
$ws.Range("E10").Value = "  -6.93%  "
$ws.Range("D12").Value = "'0.344"
$ws.Range("E12").Value = "  -2.66%  "
$ws.Range("D13").Value = "'4.85"
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("D14").Value = "2.945.14"
$ws.Range("E14").Value = "  -5.35%  "
$ws.Range("D15").Value = "69.349.06"
$ws.Range("E15").Value = "  -3.78%  "
$ws.Range("D16").Value = "'0.0000173"
$ws.Range("E16").Value = "  -5.91%  "
$ws.Range("D17").Value = "'24.84"
$ws.Range("E17").Value = "  -3.60%  "
$ws.Range("D18").Value = "2.510.85"
$ws.Range("D19").Value = "'7.81"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("D20").Value = "'11.37"
$ws.Range("E20").Value = "  -5.90%  "
$ws.Range("D21").Value = "'346.93"
$ws.Range("E21").Value = "  -7.02%  "
$ws.Range("D22").Value = "'3.94"
$ws.Range("E22").Value = "  -3.82%  "
$ws.Range("E23").Value = "  -5.49%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "'68.57"
$ws.Range("E25").Value = "  -3.16%  "
$ws.Range("D26").Value = "'3.96"
$ws.Range("E26").Value = "  -5.72%  "
$ws.Range("D27").Value = "'8.87"
$ws.Range("E27").Value = "  -7.06%  "
$ws.Range("E28").Value = "  -4.98%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").Value = "0.0₃0899"
$ws.Range("E30").Value = "  -4.88%  "
$ws.Range("D31").Value = "'7.90"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").Value = "'1.25"
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("D33").Value = "'461.82"
$ws.Range("E33").Value = "  -6.42%  "
$ws.Range("D34").Value = "'1.75"
$ws.Range("E34").Value = "  -2.13%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'0.117"
$ws.Range("E36").Value = "  +3.23%  "
$ws.Range("D37").Value = "'154.16"
$ws.Range("E37").Value = "  -5.46%  "
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("E39").Value = "  -4.05%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("E41").Value = "  -2.83%  "
$ws.Range("E42").Value = "  -2.47%  "
$ws.Range("E43").Value = "  -6.71%  "
$ws.Range("E44").Value = "  -14.19%  "
$ws.Range("E45").Value = "  -9.96%  "
$ws.Range("D46").Value = "'38.06"
$ws.Range("E46").Value = "  -2.43%  "
$ws.Range("D47").Value = "'143.40"
$ws.Range("E47").Value = "  -5.49%  "
$ws.Range("D48").Value = "'0.526"
$ws.Range("E48").Value = "  -2.95%  "
$ws.Range("D49").Value = "'3.50"
$ws.Range("E49").Value = "  -3.68%  "
$ws.Range("E50").Value = "  -4.43%  "
$ws.Range("D51").Value = "'0.0731"
$ws.Range("E51").Value = "  -1.38%  "
